# HOT-27 - Saving optimizer comparison.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "None *4" -> "Sensitive *4" (affects F5, G5, H5 since they share the same string)
$ws.Range("F5").Value = "Sensitive *4"
$ws.Range("G5").Value = "Sensitive *4"
$ws.Range("H5").Value = "Sensitive *4"

# "Sequences" -> "Sequences available"
$ws.Range("A7").Value = "Sequences available"

# "Yes *12" -> "Since 5.7"
$ws.Range("H15").Value = "Since 5.7"

# "*2 Not used internally." -> "*2 Informed by the database, but not used internally."
$ws.Range("A29").Value = "*2 Informed by the database, but not used internally."

# "*12 Since MySQL 5.7" -> "*12 "
$ws.Range("A39").Value = "*12 "

# Update the active cell selection from A8 to A40
$ws.Activate()
$ws.Range("A40").Select()
